$d = $word.ActiveDocument

# Update the title/date line.
$d.Content.Find.Execute("2023-12-15 Friday", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2023-12-16 Saturday", 2)

# Update the division problems table. Cells are addressed by (row, column)
# so that the edit is unambiguous regardless of any duplicate / reordered
# text values elsewhere in the table.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "11÷6="
$t.Cell(1, 2).Range.Text  = "28÷2="
$t.Cell(1, 3).Range.Text  = "92÷6="
$t.Cell(1, 4).Range.Text  = "62÷7="
$t.Cell(1, 5).Range.Text  = "68÷2="

$t.Cell(5, 1).Range.Text  = "27÷3="
$t.Cell(5, 2).Range.Text  = "17÷6="
$t.Cell(5, 3).Range.Text  = "67÷2="
$t.Cell(5, 4).Range.Text  = "52÷3="
$t.Cell(5, 5).Range.Text  = "78÷8="

$t.Cell(9, 1).Range.Text  = "62÷2="
$t.Cell(9, 2).Range.Text  = "31÷5="
$t.Cell(9, 3).Range.Text  = "47÷5="
$t.Cell(9, 4).Range.Text  = "47÷2="
$t.Cell(9, 5).Range.Text  = "77÷8="

$t.Cell(13, 1).Range.Text = "86÷5="
$t.Cell(13, 2).Range.Text = "80÷9="
$t.Cell(13, 3).Range.Text = "55÷5="
$t.Cell(13, 4).Range.Text = "91÷5="
$t.Cell(13, 5).Range.Text = "34÷9="

$t.Cell(17, 1).Range.Text = "11÷4="
$t.Cell(17, 2).Range.Text = "11÷6="
$t.Cell(17, 3).Range.Text = "51÷6="
$t.Cell(17, 4).Range.Text = "65÷9="
$t.Cell(17, 5).Range.Text = "12÷7="
